$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A70").Value = "4de94b55e-538e-4225-93f3-303390e81ed8"
$ws.Range("B70").Value = "champignon"
$ws.Range("C70").Value = "PSQ"
